# "Missing variable rename, ggplots"
#
# The "Data/VariableNames.xlsx" lookup sheet maps the short/renamed
# variable names (column A) used elsewhere in the analysis back to their
# original source-data column names (column B). "EducationField" was
# missing from the mapping (it's renamed to "EduField" downstream), so
# add the missing row right under "Education" (row 8) and push the rest
# of the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 9 (EntireRow insert -> everything below
# shifts down one row; the new row inherits the formatting of the row
# above it - "Education", which already uses the emphasized "Normal 2"
# style family applied to the other category-header rows in this table).
$ws.Rows.Item(9).Insert() | Out-Null

# Write column B (Original Name) before column A (Modified Name) so the
# new shared-string entries land in the same order as the source data:
# "EducationField" then "EduField".
$ws.Range("B9").Value = "EducationField"
$ws.Range("A9").Value = "EduField"

# Keep the emphasized look used by the other "category header" rows in
# this table (Education, JobLevel, PerformanceRating, StockOptionLevel):
# 12pt font, vertically centered, wrapped text.
$ws.Range("B9").Font.Size = 12
$ws.Range("B9").WrapText = $true
$ws.Range("B9").VerticalAlignment = -4108

$ws.Range("A10").Select() | Out-Null
